# ITAX.xlsx - "Db BackUp On 30-08-2021"
#
# Sheet1 previously computed a tax slab in H2:J2 (Std/PPF deduction H2,
# Net Taxable Income I2 = G2-H2, Tax J2 = 20% of I2) and had a stray
# value in G4. This edit clears that unused H2:J2 calculation (keeping
# the cell styling) and instead records the existing "Net Taxable
# Income" (column G, row 2) as a labeled value pair in row 3, columns
# I/J. The stray G4 value is removed in the process (clearing it leaves
# row 4 empty, so it simply disappears from the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old per-slab calculation in H2:J2, but keep their existing
# number formatting/styles.
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Remove the stray leftover value in G4.
$ws.Range("G4").ClearContents()

# Add a new "Net Taxable Income" label/value pair at I3:J3 (reusing the
# same shared string already used by the G1 header, and the same
# #,##0 numeric style used by G2).
$ws.Range("I3").Value = "Net Taxable Income"
$ws.Range("J3").Value = 961600
$ws.Range("J3").NumberFormat = "#,##0"

# Widen column I (now holding the longer label) and shrink column J to
# fit its new, shorter numeric content.
$ws.Columns("I").ColumnWidth = 18
$ws.Columns("J").ColumnWidth = 7.6

# Leave the active selection on the new entry.
$null = $ws.Range("I4").Select()
